# Natmi following Dr Hou advice
# Rebuild the Efna5-Epha2 LR-pair table: recompute statistics for the
# existing Sending-cluster x Target-cluster combinations and add the
# missing "ECs" sending-cluster rows (rows 2-13 total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Efna5/Epha2)
$ws.Cells.Item(2, 1).Value = "ECs"   # A2
$ws.Cells.Item(2, 2).Value = "Efna5"   # B2
$ws.Cells.Item(2, 3).Value = "Epha2"   # C2
$ws.Cells.Item(2, 4).Value = "ECs"   # D2
$ws.Cells.Item(2, 5).Value = 1   # E2
$ws.Cells.Item(2, 6).Value = 0.3333333333333333   # F2
$ws.Cells.Item(2, 7).Value = 0.09353433333333333   # G2
$ws.Cells.Item(2, 8).Value = 0.280603   # H2
$ws.Cells.Item(2, 9).Value = 0.0471360296668244   # I2
$ws.Cells.Item(2, 10).Value = 0.0471360296668244   # J2
$ws.Cells.Item(2, 11).Value = 3   # K2
$ws.Cells.Item(2, 12).Value = 1   # L2
$ws.Cells.Item(2, 13).Value = 24.244885   # M2
$ws.Cells.Item(2, 14).Value = 72.734655   # N2
$ws.Cells.Item(2, 15).Value = 0.6895205882382217   # O2
$ws.Cells.Item(2, 16).Value = 0.6895205882382218   # P2
$ws.Cells.Item(2, 17).Value = 2.267729155218333   # Q2
$ws.Cells.Item(2, 18).Value = 20.409562396965   # R2
$ws.Cells.Item(2, 19).Value = 0.03250126290308303   # S2
$ws.Cells.Item(2, 20).Value = 0.03250126290308303   # T2

# Row 3: ECs -> FAPs (Efna5/Epha2)
$ws.Cells.Item(3, 1).Value = "ECs"   # A3
$ws.Cells.Item(3, 2).Value = "Efna5"   # B3
$ws.Cells.Item(3, 3).Value = "Epha2"   # C3
$ws.Cells.Item(3, 4).Value = "FAPs"   # D3
$ws.Cells.Item(3, 5).Value = 1   # E3
$ws.Cells.Item(3, 6).Value = 0.3333333333333333   # F3
$ws.Cells.Item(3, 7).Value = 0.09353433333333333   # G3
$ws.Cells.Item(3, 8).Value = 0.280603   # H3
$ws.Cells.Item(3, 9).Value = 0.0471360296668244   # I3
$ws.Cells.Item(3, 10).Value = 0.0471360296668244   # J3
$ws.Cells.Item(3, 11).Value = 3   # K3
$ws.Cells.Item(3, 12).Value = 1   # L3
$ws.Cells.Item(3, 13).Value = 0.7651789999999999   # M3
$ws.Cells.Item(3, 14).Value = 2.295537   # N3
$ws.Cells.Item(3, 15).Value = 0.02176156637523891   # O3
$ws.Cells.Item(3, 16).Value = 0.02176156637523891   # P3
$ws.Cells.Item(3, 17).Value = 0.07157050764566666   # Q3
$ws.Cells.Item(3, 18).Value = 0.644134568811   # R3
$ws.Cells.Item(3, 19).Value = 0.00102575383825983   # S3
$ws.Cells.Item(3, 20).Value = 0.00102575383825983   # T3

# Row 4: ECs -> M2 (Efna5/Epha2)
$ws.Cells.Item(4, 1).Value = "ECs"   # A4
$ws.Cells.Item(4, 2).Value = "Efna5"   # B4
$ws.Cells.Item(4, 3).Value = "Epha2"   # C4
$ws.Cells.Item(4, 4).Value = "M2"   # D4
$ws.Cells.Item(4, 5).Value = 1   # E4
$ws.Cells.Item(4, 6).Value = 0.3333333333333333   # F4
$ws.Cells.Item(4, 7).Value = 0.09353433333333333   # G4
$ws.Cells.Item(4, 8).Value = 0.280603   # H4
$ws.Cells.Item(4, 9).Value = 0.0471360296668244   # I4
$ws.Cells.Item(4, 10).Value = 0.0471360296668244   # J4
$ws.Cells.Item(4, 11).Value = 3   # K4
$ws.Cells.Item(4, 12).Value = 1   # L4
$ws.Cells.Item(4, 13).Value = 0.6351283333333334   # M4
$ws.Cells.Item(4, 14).Value = 1.905385   # N4
$ws.Cells.Item(4, 15).Value = 0.01806294655581008   # O4
$ws.Cells.Item(4, 16).Value = 0.01806294655581008   # P4
$ws.Cells.Item(4, 17).Value = 0.05940630523944445   # Q4
$ws.Cells.Item(4, 18).Value = 0.534656747155   # R4
$ws.Cells.Item(4, 19).Value = 0.0008514155847249276   # S4
$ws.Cells.Item(4, 20).Value = 0.0008514155847249276   # T4

# Row 5: ECs -> sCs (Efna5/Epha2)
$ws.Cells.Item(5, 1).Value = "ECs"   # A5
$ws.Cells.Item(5, 2).Value = "Efna5"   # B5
$ws.Cells.Item(5, 3).Value = "Epha2"   # C5
$ws.Cells.Item(5, 4).Value = "sCs"   # D5
$ws.Cells.Item(5, 5).Value = 1   # E5
$ws.Cells.Item(5, 6).Value = 0.3333333333333333   # F5
$ws.Cells.Item(5, 7).Value = 0.09353433333333333   # G5
$ws.Cells.Item(5, 8).Value = 0.280603   # H5
$ws.Cells.Item(5, 9).Value = 0.0471360296668244   # I5
$ws.Cells.Item(5, 10).Value = 0.0471360296668244   # J5
$ws.Cells.Item(5, 11).Value = 3   # K5
$ws.Cells.Item(5, 12).Value = 1   # L5
$ws.Cells.Item(5, 13).Value = 9.516752666666667   # M5
$ws.Cells.Item(5, 14).Value = 28.550258   # N5
$ws.Cells.Item(5, 15).Value = 0.2706548988307292   # O5
$ws.Cells.Item(5, 16).Value = 0.2706548988307293   # P5
$ws.Cells.Item(5, 17).Value = 0.890143116174889   # Q5
$ws.Cells.Item(5, 18).Value = 8.011288045574   # R5
$ws.Cells.Item(5, 19).Value = 0.01275759734075661   # S5
$ws.Cells.Item(5, 20).Value = 0.01275759734075661   # T5

# Row 6: FAPs -> ECs (Efna5/Epha2)
$ws.Cells.Item(6, 1).Value = "FAPs"   # A6
$ws.Cells.Item(6, 2).Value = "Efna5"   # B6
$ws.Cells.Item(6, 3).Value = "Epha2"   # C6
$ws.Cells.Item(6, 4).Value = "ECs"   # D6
$ws.Cells.Item(6, 5).Value = 3   # E6
$ws.Cells.Item(6, 6).Value = 1   # F6
$ws.Cells.Item(6, 7).Value = 1.367717666666667   # G6
$ws.Cells.Item(6, 8).Value = 4.103153   # H6
$ws.Cells.Item(6, 9).Value = 0.6892525793933763   # I6
$ws.Cells.Item(6, 10).Value = 0.6892525793933763   # J6
$ws.Cells.Item(6, 11).Value = 3   # K6
$ws.Cells.Item(6, 12).Value = 1   # L6
$ws.Cells.Item(6, 13).Value = 24.244885   # M6
$ws.Cells.Item(6, 14).Value = 72.734655   # N6
$ws.Cells.Item(6, 15).Value = 0.6895205882382217   # O6
$ws.Cells.Item(6, 16).Value = 0.6895205882382218   # P6
$ws.Cells.Item(6, 17).Value = 33.16015754080166   # Q6
$ws.Cells.Item(6, 18).Value = 298.441417867215   # R6
$ws.Cells.Item(6, 19).Value = 0.4752538439880324   # S6
$ws.Cells.Item(6, 20).Value = 0.4752538439880325   # T6

# Row 7: FAPs -> FAPs (Efna5/Epha2)
$ws.Cells.Item(7, 1).Value = "FAPs"   # A7
$ws.Cells.Item(7, 2).Value = "Efna5"   # B7
$ws.Cells.Item(7, 3).Value = "Epha2"   # C7
$ws.Cells.Item(7, 4).Value = "FAPs"   # D7
$ws.Cells.Item(7, 5).Value = 3   # E7
$ws.Cells.Item(7, 6).Value = 1   # F7
$ws.Cells.Item(7, 7).Value = 1.367717666666667   # G7
$ws.Cells.Item(7, 8).Value = 4.103153   # H7
$ws.Cells.Item(7, 9).Value = 0.6892525793933763   # I7
$ws.Cells.Item(7, 10).Value = 0.6892525793933763   # J7
$ws.Cells.Item(7, 11).Value = 3   # K7
$ws.Cells.Item(7, 12).Value = 1   # L7
$ws.Cells.Item(7, 13).Value = 0.7651789999999999   # M7
$ws.Cells.Item(7, 14).Value = 2.295537   # N7
$ws.Cells.Item(7, 15).Value = 0.02176156637523891   # O7
$ws.Cells.Item(7, 16).Value = 0.02176156637523891   # P7
$ws.Cells.Item(7, 17).Value = 1.046548836462333   # Q7
$ws.Cells.Item(7, 18).Value = 9.418939528161   # R7
$ws.Cells.Item(7, 19).Value = 0.01499921575577358   # S7
$ws.Cells.Item(7, 20).Value = 0.01499921575577359   # T7

# Row 8: FAPs -> M2 (Efna5/Epha2)
$ws.Cells.Item(8, 1).Value = "FAPs"   # A8
$ws.Cells.Item(8, 2).Value = "Efna5"   # B8
$ws.Cells.Item(8, 3).Value = "Epha2"   # C8
$ws.Cells.Item(8, 4).Value = "M2"   # D8
$ws.Cells.Item(8, 5).Value = 3   # E8
$ws.Cells.Item(8, 6).Value = 1   # F8
$ws.Cells.Item(8, 7).Value = 1.367717666666667   # G8
$ws.Cells.Item(8, 8).Value = 4.103153   # H8
$ws.Cells.Item(8, 9).Value = 0.6892525793933763   # I8
$ws.Cells.Item(8, 10).Value = 0.6892525793933763   # J8
$ws.Cells.Item(8, 11).Value = 3   # K8
$ws.Cells.Item(8, 12).Value = 1   # L8
$ws.Cells.Item(8, 13).Value = 0.6351283333333334   # M8
$ws.Cells.Item(8, 14).Value = 1.905385   # N8
$ws.Cells.Item(8, 15).Value = 0.01806294655581008   # O8
$ws.Cells.Item(8, 16).Value = 0.01806294655581008   # P8
$ws.Cells.Item(8, 17).Value = 0.8686762421005556   # Q8
$ws.Cells.Item(8, 18).Value = 7.818086178905   # R8
$ws.Cells.Item(8, 19).Value = 0.0124499325050368   # S8
$ws.Cells.Item(8, 20).Value = 0.0124499325050368   # T8

# Row 9: FAPs -> sCs (Efna5/Epha2)
$ws.Cells.Item(9, 1).Value = "FAPs"   # A9
$ws.Cells.Item(9, 2).Value = "Efna5"   # B9
$ws.Cells.Item(9, 3).Value = "Epha2"   # C9
$ws.Cells.Item(9, 4).Value = "sCs"   # D9
$ws.Cells.Item(9, 5).Value = 3   # E9
$ws.Cells.Item(9, 6).Value = 1   # F9
$ws.Cells.Item(9, 7).Value = 1.367717666666667   # G9
$ws.Cells.Item(9, 8).Value = 4.103153   # H9
$ws.Cells.Item(9, 9).Value = 0.6892525793933763   # I9
$ws.Cells.Item(9, 10).Value = 0.6892525793933763   # J9
$ws.Cells.Item(9, 11).Value = 3   # K9
$ws.Cells.Item(9, 12).Value = 1   # L9
$ws.Cells.Item(9, 13).Value = 9.516752666666667   # M9
$ws.Cells.Item(9, 14).Value = 28.550258   # N9
$ws.Cells.Item(9, 15).Value = 0.2706548988307292   # O9
$ws.Cells.Item(9, 16).Value = 0.2706548988307293   # P9
$ws.Cells.Item(9, 17).Value = 13.01623075149711   # Q9
$ws.Cells.Item(9, 18).Value = 117.146076763474   # R9
$ws.Cells.Item(9, 19).Value = 0.1865495871445334   # S9
$ws.Cells.Item(9, 20).Value = 0.1865495871445335   # T9

# Row 10: sCs -> ECs (Efna5/Epha2)
$ws.Cells.Item(10, 1).Value = "sCs"   # A10
$ws.Cells.Item(10, 2).Value = "Efna5"   # B10
$ws.Cells.Item(10, 3).Value = "Epha2"   # C10
$ws.Cells.Item(10, 4).Value = "ECs"   # D10
$ws.Cells.Item(10, 5).Value = 3   # E10
$ws.Cells.Item(10, 6).Value = 1   # F10
$ws.Cells.Item(10, 7).Value = 0.523097   # G10
$ws.Cells.Item(10, 8).Value = 1.569291   # H10
$ws.Cells.Item(10, 9).Value = 0.2636113909397994   # I10
$ws.Cells.Item(10, 10).Value = 0.2636113909397994   # J10
$ws.Cells.Item(10, 11).Value = 3   # K10
$ws.Cells.Item(10, 12).Value = 1   # L10
$ws.Cells.Item(10, 13).Value = 24.244885   # M10
$ws.Cells.Item(10, 14).Value = 72.734655   # N10
$ws.Cells.Item(10, 15).Value = 0.6895205882382217   # O10
$ws.Cells.Item(10, 16).Value = 0.6895205882382218   # P10
$ws.Cells.Item(10, 17).Value = 12.682426608845   # Q10
$ws.Cells.Item(10, 18).Value = 114.141839479605   # R10
$ws.Cells.Item(10, 19).Value = 0.1817654813471063   # S10
$ws.Cells.Item(10, 20).Value = 0.1817654813471064   # T10

# Row 11: sCs -> FAPs (Efna5/Epha2)
$ws.Cells.Item(11, 1).Value = "sCs"   # A11
$ws.Cells.Item(11, 2).Value = "Efna5"   # B11
$ws.Cells.Item(11, 3).Value = "Epha2"   # C11
$ws.Cells.Item(11, 4).Value = "FAPs"   # D11
$ws.Cells.Item(11, 5).Value = 3   # E11
$ws.Cells.Item(11, 6).Value = 1   # F11
$ws.Cells.Item(11, 7).Value = 0.523097   # G11
$ws.Cells.Item(11, 8).Value = 1.569291   # H11
$ws.Cells.Item(11, 9).Value = 0.2636113909397994   # I11
$ws.Cells.Item(11, 10).Value = 0.2636113909397994   # J11
$ws.Cells.Item(11, 11).Value = 3   # K11
$ws.Cells.Item(11, 12).Value = 1   # L11
$ws.Cells.Item(11, 13).Value = 0.7651789999999999   # M11
$ws.Cells.Item(11, 14).Value = 2.295537   # N11
$ws.Cells.Item(11, 15).Value = 0.02176156637523891   # O11
$ws.Cells.Item(11, 16).Value = 0.02176156637523891   # P11
$ws.Cells.Item(11, 17).Value = 0.400262839363   # Q11
$ws.Cells.Item(11, 18).Value = 3.602365554267   # R11
$ws.Cells.Item(11, 19).Value = 0.005736596781205498   # S11
$ws.Cells.Item(11, 20).Value = 0.005736596781205499   # T11

# Row 12: sCs -> M2 (Efna5/Epha2)
$ws.Cells.Item(12, 1).Value = "sCs"   # A12
$ws.Cells.Item(12, 2).Value = "Efna5"   # B12
$ws.Cells.Item(12, 3).Value = "Epha2"   # C12
$ws.Cells.Item(12, 4).Value = "M2"   # D12
$ws.Cells.Item(12, 5).Value = 3   # E12
$ws.Cells.Item(12, 6).Value = 1   # F12
$ws.Cells.Item(12, 7).Value = 0.523097   # G12
$ws.Cells.Item(12, 8).Value = 1.569291   # H12
$ws.Cells.Item(12, 9).Value = 0.2636113909397994   # I12
$ws.Cells.Item(12, 10).Value = 0.2636113909397994   # J12
$ws.Cells.Item(12, 11).Value = 3   # K12
$ws.Cells.Item(12, 12).Value = 1   # L12
$ws.Cells.Item(12, 13).Value = 0.6351283333333334   # M12
$ws.Cells.Item(12, 14).Value = 1.905385   # N12
$ws.Cells.Item(12, 15).Value = 0.01806294655581008   # O12
$ws.Cells.Item(12, 16).Value = 0.01806294655581008   # P12
$ws.Cells.Item(12, 17).Value = 0.3322337257816667   # Q12
$ws.Cells.Item(12, 18).Value = 2.990103532035   # R12
$ws.Cells.Item(12, 19).Value = 0.004761598466048354   # S12
$ws.Cells.Item(12, 20).Value = 0.004761598466048354   # T12

# Row 13: sCs -> sCs (Efna5/Epha2)
$ws.Cells.Item(13, 1).Value = "sCs"   # A13
$ws.Cells.Item(13, 2).Value = "Efna5"   # B13
$ws.Cells.Item(13, 3).Value = "Epha2"   # C13
$ws.Cells.Item(13, 4).Value = "sCs"   # D13
$ws.Cells.Item(13, 5).Value = 3   # E13
$ws.Cells.Item(13, 6).Value = 1   # F13
$ws.Cells.Item(13, 7).Value = 0.523097   # G13
$ws.Cells.Item(13, 8).Value = 1.569291   # H13
$ws.Cells.Item(13, 9).Value = 0.2636113909397994   # I13
$ws.Cells.Item(13, 10).Value = 0.2636113909397994   # J13
$ws.Cells.Item(13, 11).Value = 3   # K13
$ws.Cells.Item(13, 12).Value = 1   # L13
$ws.Cells.Item(13, 13).Value = 9.516752666666667   # M13
$ws.Cells.Item(13, 14).Value = 28.550258   # N13
$ws.Cells.Item(13, 15).Value = 0.2706548988307292   # O13
$ws.Cells.Item(13, 16).Value = 0.2706548988307293   # P13
$ws.Cells.Item(13, 17).Value = 4.978184769675334   # Q13
$ws.Cells.Item(13, 18).Value = 44.803662927078   # R13
$ws.Cells.Item(13, 19).Value = 0.07134771434543923   # S13
$ws.Cells.Item(13, 20).Value = 0.07134771434543924   # T13

